$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels in row 1 (columns AE:AI) ---
# The shared-string usage-label columns got reshuffled:
#   AE1: tkm-SZMUsage -> pkmUsage
#   AF1: tkm-N3Usage  -> tkm-SZMUsage
#   AG1: tkm-N2Usage  -> tkm-N1Usage
#   AH1: pkmUsage     -> tkm-N2Usage
#   AI1: tkm-N1Usage  -> tkm-N3Usage
$ws.Range("AE1").Value = "pkmUsage"
$ws.Range("AF1").Value = "tkm-SZMUsage"
$ws.Range("AG1").Value = "tkm-N1Usage"
$ws.Range("AH1").Value = "tkm-N2Usage"
$ws.Range("AI1").Value = "tkm-N3Usage"

# --- Update data row 2 values to match the reshuffled/recomputed figures ---
$ws.Range("B2").Value = 42.50739806315688
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 612.0528792959246
$ws.Range("H2").Value = 741.7852278688524
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = 514.1244186085767
$ws.Range("N2").Value = 514.1244186085767
$ws.Range("R2").Value = 288.4396604831751

$ws.Range("AE2").Value = 858
$ws.Range("AF2").Value = 414.5
$ws.Range("AG2").Value = 7.5
$ws.Range("AH2").Value = 24.2
$ws.Range("AI2").Value = 130.3
